# Apply the Oct 22 2023 cryptos-list refresh (prices / 1h volumes + two rank swaps).
# Values are plain text in the sheet (e.g. "29.821.56", "  -0.21%  "), so we assign
# strings throughout. Numeric-looking text (e.g. "215.17", "1.00") gets a leading
# apostrophe, same as typing it by hand in Excel, so it is stored as text (quotePrefix)
# instead of being silently parsed into a float and losing formatting/precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '29.821.56'
$ws.Range("E2").Value = '  -0.21%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.638.14'
$ws.Range("E3").Value = '  +0.89%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.50%  '

# Row 5: BNB
$ws.Range("D5").Value = '''215.17'
$ws.Range("E5").Value = '  +0.41%  '

# Row 6: XRP
$ws.Range("E6").Value = '  -0.42%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.49%  '

# Row 8: Solana
$ws.Range("D8").Value = '''28.92'
$ws.Range("E8").Value = '  -3.65%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  +0.76%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  -0.27%  '

# Row 11: TRON
$ws.Range("D11").Value = '''0.0899'

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '1.873.06'
$ws.Range("E12").Value = '  +0.95%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '1.636.23'
$ws.Range("E13").Value = '  +1.09%  '

# Row 14: Polygon
$ws.Range("D14").Value = '''0.591'
$ws.Range("E14").Value = '  +3.93%  '

# Row 15: Chainlink
$ws.Range("D15").Value = '''9.52'
$ws.Range("E15").Value = '  +7.78%  '

# Row 16: Polkadot
$ws.Range("E16").Value = '  -0.38%  '

# Row 17: WrappedBTC
$ws.Range("D17").Value = '29.844.18'
$ws.Range("E17").Value = '  -0.29%  '

# Row 18: Litecoin
$ws.Range("D18").Value = '''64.20'
$ws.Range("E18").Value = '  -0.64%  '

# Row 19: BitcoinCash
$ws.Range("D19").Value = '''238.32'
$ws.Range("E19").Value = '  -2.33%  '

# Row 20: ShibaInu
$ws.Range("E20").Value = '  -0.22%  '

# Row 21: Dai
$ws.Range("E21").Value = '  +0.40%  '

# Row 22: Avalanche
$ws.Range("E22").Value = '  +2.97%  '

# Row 23: Uniswap
$ws.Range("E23").Value = '  +0.23%  '

# Row 24: Toncoin
$ws.Range("E24").Value = '  +2.13%  '

# Row 25: Monero
$ws.Range("E25").Value = '  +0.22%  '

# Row 26: EthereumClassic
$ws.Range("D26").Value = '''15.57'
$ws.Range("E26").Value = '  -0.66%  '

# Row 27: Stellar
$ws.Range("E27").Value = '  -1.12%  '

# Row 28: Cosmos
$ws.Range("E28").Value = '  +0.18%  '

# Row 29: BinanceUSD
$ws.Range("E29").Value = '  +0.44%  '

# Row 30: Hedera
$ws.Range("E30").Value = '  +1.29%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '  -0.66%  '

# Row 32: Filecoin
$ws.Range("E32").Value = '  +1.24%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = '  -0.96%  '

# Row 34: Maker
$ws.Range("D34").Value = '1.418.33'
$ws.Range("E34").Value = '  -0.61%  '

# Row 35: LidoDAOToken
$ws.Range("E35").Value = '  +1.95%  '

# Row 36: TrustWalletToken
$ws.Range("E36").Value = '  -1.25%  '

# Row 37: MXToken
$ws.Range("D37").Value = '''2.71'
$ws.Range("E37").Value = '  -5.68%  '

# Row 38: VeChain
$ws.Range("E38").Value = '  +2.18%  '

# Row 39: HuobiToken
$ws.Range("E39").Value = '  +0.32%  '

# Row 40: Aave
$ws.Range("E40").Value = '  +10.54%  '

# Row 41: ImmutableX
$ws.Range("E41").Value = '  +1.06%  '

# Row 42: row 42 (ARBITRUM -> Kaspa, rank swap with row 43)
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.0501'
$ws.Range("E42").Value = '  -0.87%  '

# Row 43: row 43 (Kaspa -> ARBITRUM, rank swap with row 42)
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '''0.832'
$ws.Range("E43").Value = '  -0.16%  '

# Row 44: RenderToken
$ws.Range("E44").Value = '  -3.16%  '

# Row 45: PaxDollar
$ws.Range("E45").Value = '  +0.48%  '

# Row 46: WEMIXToken
$ws.Range("D46").Value = '''1.00'
$ws.Range("E46").Value = '  -2.14%  '

# Row 47: row 47 (RocketPoolETH -> BitcoinSV, rank swap with row 48)
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").Value = '''49.95'
$ws.Range("E47").Value = '  -7.53%  '

# Row 48: row 48 (BitcoinSV -> RocketPoolETH, rank swap with row 47)
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.781.79'
$ws.Range("E48").Value = '  +1.00%  '

# Row 49: FraxShare
$ws.Range("D49").Value = '''5.33'
$ws.Range("E49").Value = '  -1.37%  '

# Row 50: Quant
$ws.Range("D50").Value = '''93.26'
$ws.Range("E50").Value = '  +5.42%  '

# Row 51: BabyDogeCoin
$ws.Range("D51").Value = '0.0₆0109'
$ws.Range("E51").Value = '  +1.20%  '
